$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handback transform failed" ---
# (every cell that used to show "Ready for handoff" needs updating so the
#  text change shows up everywhere it was referenced)
$wsOverview.Range("E2").Value = "Handback transform failed"
$wsOverview.Range("F2").Value = "Handback transform failed"
$wsZhCn.Range("C2").Value = "Handback transform failed"
$wsDeDe.Range("C2").Value = "Handback transform failed"

# --- New "Error Detail" messages for the handback-priority mismatch ---
$wsZhCn.Range("O2").Value = "The handback priority in file path mt\032a9fc0-b57a-43c8-bcef-3e8dfd2ba006.4f450a64dfbe2dab1d19b9d182a24e07c825aaca.zh-cn.xlf is not match with handoff type ht."
$wsDeDe.Range("O2").Value = "The handback priority in file path mt\032a9fc0-b57a-43c8-bcef-3e8dfd2ba006.4f450a64dfbe2dab1d19b9d182a24e07c825aaca.de-de.xlf is not match with handoff type ht."

# --- Column width adjustments ---
# Overview: zh-cn / de-de status columns widened (17.216 -> 24.7427 chars)
$wsOverview.Columns.Item(5).ColumnWidth = 23.833333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 23.833333333333336

# zh-cn / de-de: Status column widened (17.216 -> 24.7427 chars)
$wsZhCn.Columns.Item(3).ColumnWidth = 23.833333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 23.833333333333336

# zh-cn / de-de: Error Detail column widened (13.747 -> 40 chars)
$wsZhCn.Columns.Item(15).ColumnWidth = 39.16666666666667
$wsDeDe.Columns.Item(15).ColumnWidth = 39.16666666666667
